{"js": "// Remove the unused \"Abstract Title\" paragraph style and change the\n// \"Abstract\" style's space-before from 5pt (100 twips) to 15pt (300 twips).\n\nconst styles = context.document.getStyles();\n\n// Delete the \"AbstractTitle\" custom style (display name \"Abstract Title\").\nconst abstractTitleStyle = styles.getByNameOrNullObject(\"Abstract Title\");\nawait context.sync();\n\nif (!abstractTitleStyle.isNullObject) {\n  abstractTitleStyle.delete();\n  await context.sync();\n}\n\n// Update the \"Abstract\" style's paragraph spacing: before 100 -> 300 twips\n// (5pt -> 15pt). spaceAfter (300 twips = 15pt) stays unchanged.\nconst abstractStyle = styles.getByName(\"Abstract\");\nabstractStyle.paragraphFormat.spaceBefore = 15;\nawait context.sync();\n", "ps1": "# Remove the unused \"Abstract Title\" paragraph style and change the\n# \"Abstract\" style's space-before from 5pt (100 twips) to 15pt (300 twips).\n\n$d = $word.ActiveDocument\n\n# Delete the \"AbstractTitle\" custom style (display name \"Abstract Title\"),\n# if present.\nforeach ($s in $d.Styles) {\n    if ($s.NameLocal -eq \"Abstract Title\") {\n        $s.Delete()\n        break\n    }\n}\n\n# Update the \"Abstract\" style's paragraph spacing before: 5pt -> 15pt\n# (100 -> 300 twips). SpaceAfter (300 twips = 15pt) stays unchanged.\n$abstract = $d.Styles.Item(\"Abstract\")\n$abstract.ParagraphFormat.SpaceBefore = 15\n\n"}
